# Convert the RANDBETWEEN shared-formula results in column I (rows 2-14)
# of the "Reporte de transacciones" sheet into static values, as happens
# when you copy a range and paste it back as values-only (breaking the
# link to the volatile RANDBETWEEN formula while keeping the last
# calculated numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value  = 63389805
$ws.Range("I3").Value  = 35566129
$ws.Range("I4").Value  = 92217107
$ws.Range("I5").Value  = 10830860
$ws.Range("I6").Value  = 24364532
$ws.Range("I7").Value  = 8202621
$ws.Range("I8").Value  = 30726945
$ws.Range("I9").Value  = 73110557
$ws.Range("I10").Value = 56156614
$ws.Range("I11").Value = 57263095
$ws.Range("I12").Value = 2778550
$ws.Range("I13").Value = 1477475
$ws.Range("I14").Value = 54090044
